$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Header numeric updates -------------------------------------------------
$ws.Range("E11").Value = 260000
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 5

# --- Swap the worker shown in the first block of rows (16-19): was
#     DIANILYS MELISSA VILLADIEGO CASTRO / CC 1043647690, now NATALY CORTES
#     SIMANCA / CC 1007978849 - the same worker already used further down
#     the table. Using Replace (instead of re-typing the values) folds the
#     two text blocks into one shared entry exactly like Excel would when
#     the resulting text becomes identical. ----------------------------------
$ws.Cells.Replace("1043647690", "1007978849")
$ws.Cells.Replace("DIANILYS MELISSA VILLADIEGO CASTRO", "NATALY CORTES SIMANCA")

# Periods for that first block move from 2507/2506/2505/2504 (descending) to
# 2504/2505/2506/2507 (ascending), matching the rest of the refreshed table.
$ws.Cells.Item(16, 5).Value = "2504"
$ws.Cells.Item(17, 5).Value = "2505"
$ws.Cells.Item(18, 5).Value = "2506"
$ws.Cells.Item(19, 5).Value = "2507"

# Their salary figures drop from 56940/1423500 to 52000/1300000, matching the
# NATALY rows already below them.
for ($r = 16; $r -le 19; $r++) {
    $ws.Cells.Item($r, 6).Value = 52000
    $ws.Cells.Item($r, 7).Value = 1300000
}

# --- Row 20 becomes the new (and last) data row for NATALY, period 2508.
#     Grab the "last row" bottom-border formatting from row 25 before that
#     row is swept away by the deletion below. ------------------------------
$ws.Range("B25:J25").Copy() | Out-Null
$ws.Range("B20:J20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Cells.Item(20, 5).Value = "2508"

# --- Drop the now-redundant rows (old rows 21-25); this also slides the
#     signature block (old rows 30-31) up to rows 25-26. ---------------------
$ws.Rows("21:25").Delete()

# --- Column D no longer needs to fit the long former name; narrow it to the
#     new best-fit width. -----------------------------------------------------
$ws.Columns("D").ColumnWidth = 24.54296875
